$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; this shifts the existing
# yearly data (previously in D:J) one column to the right (E:K)
# and creates a brand-new, blank column D for the latest fiscal year.
$ws.Columns("D:D").Insert()

# The newly inserted column D has no number formatting; copy the
# formatting that column E (the old column D) already carries so the
# new column matches the rest of the data (date format / number format).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New cell values (columns D through L) for every row that changed.
# Most rows simply gained a new "D" value and shifted the rest one
# column right, but several historical figures were also restated,
# so the full row contents are listed explicitly.
$rowData = @{
    7 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908, $null)
    8 = @(518100, 426000, 400900, 336800, 312900, 265900, 244300, 230100, $null)
    9 = @(379800, 310100, 294600, 245500, 229200, 198200, 180300, 168600, $null)
    10 = @(138300, 115900, 106300, 91300, 83700, 67700, 64100, 61500, $null)
    12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", $null)
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    14 = @("NA", -2500, 0, 100, "NA", 0, 0, 15500, $null)
    15 = @(8600, 6700, 6600, 4700, 3800, 2200, 2500, 3200, $null)
    17 = @(493800, 399600, 385500, 320800, 294900, 250500, 228700, 217200, $null)
    18 = @(24300, 26400, 15500, 16000, 18100, 15500, 15700, 12900, $null)
    20 = @(2600, 300, 3000, 0, 0, 200, 200, 2300, $null)
    21 = @(35500, 33300, 25100, 20800, 21900, 17800, 18400, 18700, $null)
    22 = @(5000, 4500, 2300, 800, 700, 700, 1700, 2500, $null)
    23 = @(21900, 22200, 16200, 15300, 17400, 15000, 14100, 12700, $null)
    24 = @(4500, 8600, 4100, 3900, 5400, 3800, 4800, 4200, $null)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    26 = @(17400, 13500, 12100, 11400, 12000, 11200, 9300, 8400, $null)
    27 = @(17400, 13500, 12100, 11400, 12000, 11200, 9300, 8400, $null)
    28 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    29 = @(100, -1200, 100, 300, 300, 8000, -1700, -10400, $null)
    30 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    32 = @(-2600, -300, -3000, 0, 0, -200, -200, -2300, $null)
    33 = @(17500, 12400, 12200, 11600, 12200, 19100, 7600, -2000, $null)
    34 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    35 = @(17500, 12400, 12200, 11600, 12200, 19100, 7600, -2000, $null)
    38 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908, $null)
    41 = @(70400, 53800, 8000, 4100, 13400, 15600, 1700, 2000, $null)
    42 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    43 = @(109700, 94100, 117700, 86300, 69700, 62200, 72300, 74200, $null)
    44 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    45 = @(5400, 7800, 5300, 12200, 14300, 13700, 13800, 18100, $null)
    46 = @(185500, 155700, 131000, 93900, 97400, 91500, 87800, 89100, $null)
    47 = @("NA", "NA", 900, 900, 900, 900, "NA", "NA", $null)
    48 = @(10700, 15000, 6600, 8600, 7700, 5300, 2500, 2300, $null)
    49 = @(159200, 106900, 88100, 79200, 74600, 68800, 56900, 66800, $null)
    50 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    51 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    52 = @(0, 1600, 3400, 3200, 300, 100, 2600, 4600, $null)
    53 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    54 = @(355400, 271700, 229900, 185800, 180800, 163900, 149900, 154700, $null)
    57 = @(12200, 7400, 4500, 4700, 4000, 4600, 4100, 5300, $null)
    58 = @(100, 3100, 2500, 1100, 1000, 0, 200, 6600, $null)
    59 = @(49200, 46400, 42600, 36300, 38300, 42100, 34900, 60800, $null)
    60 = @(61500, 55100, 49600, 42200, 43200, 46600, 39200, 43300, $null)
    61 = @(17200, 39900, 22500, 1900, 2700, 0, 16300, 25000, $null)
    62 = @(1100, 400, "NA", 6800, 7000, 3400, "NA", 0, $null)
    63 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    64 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    65 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    66 = @(79900, 95400, 72100, 44100, 52800, 50100, 55400, 68300, $null)
    68 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    69 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    70 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    71 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    72 = @(97800, 80300, 65500, 54600, 43000, 30800, 11600, 4000, $null)
    73 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    74 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    75 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    76 = @(275500, 176300, 157800, 141700, 128000, 113900, 94400, 86400, $null)
    77 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    80 = @(43465, 43100, 42735, 42369, 42004, 41639, 41274, 40908, $null)
    81 = @(17500, 12400, 12200, 11600, 12200, 19100, 7600, -2000, $null)
    83 = @(8600, 6700, 6600, 4700, 3800, 2200, 2500, 3600, $null)
    84 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    85 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    86 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    87 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    88 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    89 = @(33200, 52800, -700, 4100, 7000, 27400, 15400, 15900, $null)
    91 = @(-5300, -3600, -1700, -2200, -5300, -900, -1100, -600, $null)
    92 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    93 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    94 = @(-67800, -24300, -21700, -10700, -12500, 2900, -600, -1100, $null)
    96 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    97 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    98 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    99 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    100 = @(51200, 17200, 26400, -2600, 3300, -16500, -15100, -13700, $null)
    101 = @(0, 0, 0, 0, 0, 0, 0, 0, $null)
    102 = @(16700, 45700, 3900, -9300, -2200, 13800, -300, 1200, $null)
}

foreach ($rowNum in $rowData.Keys) {
    $vals = $rowData[$rowNum]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $v = $vals[$i]
        if ($v -ne $null) {
            $ws.Cells.Item($rowNum, 4 + $i).Value = $v
        }
    }
}
